$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(-0.02857142857142796, 0.7571428571428517, 1.618571428571423, 1.272230886502691, 1.319923407701294)
    3  = @(-0.03076923076922793, 0.8923076923076908, 2.001538461538451, 1.414757386104929, 1.472177873868399)
    4  = @(0.0250000000000045, 0.8416666666666673, 1.390833333333325, 1.179335971355629, 1.231499455578808)
    5  = @(-0.07272727272726601, 0.7818181818181839, 1.32727272727272, 1.152073229995698, 1.205894613215355)
    6  = @(-0.139999999999992, 0.9199999999999988, 1.737999999999984, 1.318332279814154, 1.381786283523367)
    7  = @(0.1333333333333415, 0.7111111111111149, 1.12888888888889, 1.062491830033949, 1.118033988749894)
    8  = @(-0.3499999999999896, 0.5166666666666609, 0.3583333333333243, 0.5986094998689248, 0.5319774431308135)
    9  = @(-0.2999999999999924, 0.7000000000000076, 0.5766666666666664, 0.7593857166596343, 0.854400374531757)
    10 = @(0.6000000000000227, 0.6000000000000227, 0.3600000000000273, 0.6000000000000227, $null)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
    $ws.Cells.Item($row, 5).Value = $rowValues[3]
    if ($null -ne $rowValues[4]) {
        $ws.Cells.Item($row, 6).Value = $rowValues[4]
    }
}
